$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A3 value (shared string "test4@test.com" -> "t@t.com")
$ws.Range("A3").Value = "t@t.com"

# Update the active selection on the sheet to D4
$ws.Range("D4").Select()
